# Apply QuantitativeEvaluation updates to the "QuantitativeMetrics" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("QuantitativeMetrics")

# Row 6 - "Runtime without error": now flagged "no" with an explanatory note.
$ws.Range("B6").Value = "no"
$ws.Range("C6").Value = "Missing initial redirect"

# Row 7 - "Assertion validity": value + note cleared out.
$ws.Range("B7").ClearContents()
$ws.Range("C7").ClearContents()

# Row 12 - "Code BLEU": refreshed score + detail breakdown.
$ws.Range("B12").Value = 0.2740119028029336
$ws.Range("C12").Value = "{'codebleu': 0.27401190280293364, 'ngram_match_score': 0.12114595985200712, 'weighted_ngram_match_score': 0.138099238086726, 'syntax_match_score': 0.5897435897435898, 'dataflow_match_score': 0.24705882352941178}"

# Update the sheet's active selection to B7, matching the last reviewed cell.
$ws.Range("B7").Select()
